# Add a new "New Vocab" / "newvocab" term list with three terms
# (Term 1, Term 2, Term 3) to the termData table, each loaded with
# loadVersion=4 / loadAction="create" - mirrors the data a user would
# type directly into the Excel table for TermListVocab.present_in_version?
# & .not_yet_loaded? spec fixtures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$table = $ws.ListObjects.Item("termData")

$terms = @("Term 1", "Term 2", "Term 3")

foreach ($term in $terms) {
    $row = $table.ListRows.Add()
    $r = $row.Range()
    $rowNum = $r.Row()

    $r.Cells.Item(1, 1).Value = 4
    $r.Cells.Item(1, 2).Value = "create"
    $r.Cells.Item(1, 4).Value = "New Vocab"
    $r.Cells.Item(1, 5).Value = "newvocab"
    $r.Cells.Item(1, 6).Value = $term
    $r.Cells.Item(1, 12).Value = $term

    $idCell = $ws.Cells.Item($rowNum, 3)
    $idCell.Formula = '=_xlfn.TEXTJOIN(" ",TRUE,E' + $rowNum + ',L' + $rowNum + ')'

    $sortCell = $ws.Cells.Item($rowNum, 11)
    $sortCell.Formula = '=_xlfn.TEXTJOIN(" ",TRUE,C' + $rowNum + ',A' + $rowNum + ')'
    $sortCell.NumberFormat = "General"
}

# Select the last term's origterm cell like the author's saved view.
$lastRow = $table.ListRows.Count() + 1
$ws.Range("L" + ($lastRow - 1) + ":L" + $lastRow).Select()
$ws.Application.ActiveWindow.ScrollRow = 27

Write-Host "Added" $terms.Count "rows; table now spans" $table.Range().Address()
